$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($addr in @("D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D17", "D20", "D21", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D50", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.704.79'
$ws.Range('E2').Value = '  +0.12%  '

$ws.Range('D3').Value = '1.904.23'
$ws.Range('E3').Value = '  +0.62%  '

$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.17%  '

$ws.Range('D5').Value = '312.43'
$ws.Range('E5').Value = '  +0.30%  '

$ws.Range('E6').Value = '  -0.10%  '

$ws.Range('D7').Value = '0.5203'
$ws.Range('E7').Value = '  +5.76%  '

$ws.Range('D8').Value = '0.3781'
$ws.Range('E8').Value = '  -0.41%  '

$ws.Range('D9').Value = '0.07256'
$ws.Range('E9').Value = '  -0.77%  '

$ws.Range('D10').Value = '21.19'
$ws.Range('E10').Value = '  +3.13%  '

$ws.Range('D11').Value = '0.9010'
$ws.Range('E11').Value = '  -0.89%  '

$ws.Range('D12').Value = '0.07664'
$ws.Range('E12').Value = '  +0.20%  '

$ws.Range('D13').Value = '1.907.02'
$ws.Range('E13').Value = '  +0.61%  '

$ws.Range('D14').Value = '5.441'
$ws.Range('E14').Value = '  -0.41%  '

$ws.Range('D15').Value = '91.94'
$ws.Range('E15').Value = '  +0.90%  '

$ws.Range('E16').Value = '  -0.20%  '

$ws.Range('D17').Value = '0.000008703'
$ws.Range('E17').Value = '  -0.50%  '

$ws.Range('E18').Value = '  -0.07%  '

$ws.Range('D19').Value = '27.757.31'
$ws.Range('E19').Value = '  +0.03%  '

$ws.Range('D20').Value = '14.51'
$ws.Range('E20').Value = '  +0.50%  '

$ws.Range('D21').Value = '5.140'
$ws.Range('E21').Value = '  +0.41%  '

$ws.Range('D22').Value = '2.152.41'
$ws.Range('E22').Value = '  +0.01%  '

$ws.Range('E23').Value = '  +0.91%  '

$ws.Range('D24').Value = '6.632'
$ws.Range('E24').Value = '  +0.08%  '

$ws.Range('D25').Value = '153.79'
$ws.Range('E25').Value = '  -0.09%  '

$ws.Range('D26').Value = '1.872'
$ws.Range('E26').Value = '  +0.85%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '18.31'
$ws.Range('E27').Value = '  -0.22%  '

$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '2.161'
$ws.Range('E28').Value = '  -0.13%  '

$ws.Range('D29').Value = '114.60'
$ws.Range('E29').Value = '  -0.45%  '

$ws.Range('D30').Value = '4.856'
$ws.Range('E30').Value = '  -0.23%  '

$ws.Range('D31').Value = '0.09042'
$ws.Range('E31').Value = '  +1.14%  '

$ws.Range('D32').Value = '3.188'
$ws.Range('E32').Value = '  -0.45%  '

$ws.Range('D33').Value = '4.833'
$ws.Range('E33').Value = '  +4.68%  '

$ws.Range('E34').Value = '  +0.69%  '

$ws.Range('D35').Value = '0.7824'
$ws.Range('E35').Value = '  +2.35%  '

$ws.Range('D36').Value = '0.02088'
$ws.Range('E36').Value = '  +2.52%  '

$ws.Range('D37').Value = '2.596'
$ws.Range('E37').Value = '  +1.74%  '

$ws.Range('D38').Value = '3.070'
$ws.Range('E38').Value = '  +2.85%  '

$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '1.093'
$ws.Range('E39').Value = '  -0.39%  '

$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '0.5565'
$ws.Range('E40').Value = '  +2.06%  '

$ws.Range('D41').Value = '0.05285'
$ws.Range('E41').Value = '  +0.03%  '

$ws.Range('D42').Value = '6.723'
$ws.Range('E42').Value = '  -2.34%  '

$ws.Range('D43').Value = '114.89'
$ws.Range('E43').Value = '  +1.33%  '

$ws.Range('D44').Value = '8.529'
$ws.Range('E44').Value = '  +0.09%  '

$ws.Range('D45').Value = '0.1520'
$ws.Range('E45').Value = '  -0.02%  '

$ws.Range('D46').Value = '0.4811'
$ws.Range('E46').Value = '  +0.57%  '

$ws.Range('D47').Value = '10.47'
$ws.Range('E47').Value = '  -0.87%  '

$ws.Range('E48').Value = '  -0.11%  '

$ws.Range('E49').Value = '  -0.94%  '

$ws.Range('D50').Value = '66.82'
$ws.Range('E50').Value = '  -0.74%  '

$ws.Range('D51').Value = '0.05991'
$ws.Range('E51').Value = '  -1.17%  '
